# Updates the cryptocurrency price/volume table on Sheet1 with refreshed
# quote data (GitHub Actions scheduled symbol-list refresh).
#
# Columns: D = Price, E = Volume(1h) change %.
# Values are stored as text in the workbook (not numbers), so each value
# is written with a leading apostrophe to force Excel to keep it as text
# instead of auto-converting the numeric-looking string to a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - BNB
$ws.Range("D2").Value = "'314.21"
$ws.Range("E2").Value = "'2.22%"

# Row 3 - OKB
$ws.Range("D3").Value = "'40.39"
$ws.Range("E3").Value = "'-1.57%"

# Row 4 - HuobiToken
$ws.Range("D4").Value = "'5.164"
$ws.Range("E4").Value = "'-1.71%"

# Row 5 - Cronos
$ws.Range("D5").Value = "'0.07584"
$ws.Range("E5").Value = "'-0.99%"

# Row 6 - GateToken
$ws.Range("D6").Value = "'4.320"
$ws.Range("E6").Value = "'0.12%"

# Row 7 - FTXToken
$ws.Range("D7").Value = "'1.670"
$ws.Range("E7").Value = "'2.74%"

# Row 8 - MXToken
$ws.Range("D8").Value = "'0.9252"
$ws.Range("E8").Value = "'0.75%"

# Row 9 - BTSEToken (price unchanged)
$ws.Range("E9").Value = "'-0.82%"

# Row 10 - LiechtensteinCryptoassetsExchange
$ws.Range("D10").Value = "'0.1201"
$ws.Range("E10").Value = "'-4.33%"

# Row 11 - WazirX
$ws.Range("D11").Value = "'0.1815"
$ws.Range("E11").Value = "'-1.03%"

# Row 12 - MandalaExchangeToken
$ws.Range("D12").Value = "'0.09043"
$ws.Range("E12").Value = "'-1.71%"

# Row 13 - BitrueCoin
$ws.Range("D13").Value = "'0.04108"
$ws.Range("E13").Value = "'-3.64%"

# Row 14 - BitMartToken (price unchanged)
$ws.Range("E14").Value = "'0.01%"

# Row 15 - BitForexToken
$ws.Range("D15").Value = "'0.001284"
$ws.Range("E15").Value = "'1.63%"

# Row 16 - TigerCash
$ws.Range("D16").Value = "'0.006004"
$ws.Range("E16").Value = "'4.50%"

# Row 18 - LEO
$ws.Range("D18").Value = "'3.351"
$ws.Range("E18").Value = "'-0.05%"

# Row 20 - MCDex
$ws.Range("D20").Value = "'7.627"
$ws.Range("E20").Value = "'2.93%"

# Row 21 - ProBitToken
$ws.Range("D21").Value = "'0.1352"
$ws.Range("E21").Value = "'-2.35%"

# Row 22 - ZBToken (price unchanged)
$ws.Range("E22").Value = "'-2.87%"

# Row 23 - CoinExToken
$ws.Range("D23").Value = "'0.04028"
$ws.Range("E23").Value = "'-1.14%"

# Row 24 - BitKan
$ws.Range("D24").Value = "'0.001273"
$ws.Range("E24").Value = "'0.76%"

# Row 25 - HotbitToken (price unchanged)
$ws.Range("E25").Value = "'-7.43%"

# Row 26 - NitroEx (price unchanged)
$ws.Range("E26").Value = "'-0.23%"

# Row 38 - One
$ws.Range("D38").Value = "'0.02420"
$ws.Range("E38").Value = "'-1.89%"

# Row 39 - IDEX
$ws.Range("D39").Value = "'0.05151"
$ws.Range("E39").Value = "'-2.48%"

# Row 40 - KickToken
$ws.Range("D40").Value = "'0.007704"
$ws.Range("E40").Value = "'-1.83%"

# Row 41 - BKEXToken
$ws.Range("D41").Value = "'0.1301"
$ws.Range("E41").Value = "'-0.95%"

# Row 42 - Dexo
$ws.Range("D42").Value = "'0.007622"
$ws.Range("E42").Value = "'11.77%"

# Row 43 - CEJI
$ws.Range("D43").Value = "'0.003303"
$ws.Range("E43").Value = "'72.54%"

# Row 44 - LocalTraders
$ws.Range("D44").Value = "'0.008260"
$ws.Range("E44").Value = "'7.18%"

# Row 45 - PooCoin
$ws.Range("D45").Value = "'0.3103"
$ws.Range("E45").Value = "'1.67%"

# Row 46 - CoinLion
$ws.Range("D46").Value = "'0.00006592"
$ws.Range("E46").Value = "'-2.04%"

# Row 47 - Kangarootoken
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'-0.17%"

# Row 48 - BOLO
$ws.Range("D48").Value = "'0.2677"
$ws.Range("E48").Value = "'57.61%"

# Row 49 - CoinbaseStockToken
$ws.Range("D49").Value = "'0.004203"
$ws.Range("E49").Value = "'2.52%"

# Row 50 - CryptobidCoin
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'-0.17%"

# Row 51 - SpecialPowerGold
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'-0.17%"
